{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Reminder: Installation -> Safety -> Tool Position...\" paragraph\n// that the two new journal entries are inserted directly after.\nconst anchor = paras.items.find(p =>\n  p.text.indexOf(\"Reminder: Installation\") !== -1 &&\n  p.text.indexOf(\"Tool Position\") !== -1\n);\n\nif (!anchor) {\n  throw new Error(\"Could not find the anchor paragraph to insert after.\");\n}\n\n// Insert in reverse order, each \"After\" the anchor, so the final order is:\n//   ... Reminder: Installation -> Safety -> Tool Position ...\n//   Reminder: set the Home position\n//   Software Safety Settings. Password is plastic\nanchor.insertParagraph(\"Software Safety Settings. Password is plastic\", Word.InsertLocation.after);\nanchor.insertParagraph(\"Reminder: set the Home position\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Reminder: Installation -> Safety -> Tool Position...\" paragraph;\n# the two new journal entries are inserted directly after it.\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Reminder: Installation*Tool Position*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find the anchor paragraph to insert after.\"\n}\n\n$anchor = $d.Paragraphs.Item($targetIndex)\n$anchor.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($targetIndex + 1).Range.Text = \"Reminder: set the Home position\"\n\n$second = $d.Paragraphs.Item($targetIndex + 1)\n$second.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($targetIndex + 2).Range.Text = \"Software Safety Settings. Password is plastic\"\n"}
